$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 996.8
$ws.Range("I33").Value = 1065.2858
$ws.Range("J33").Value = 38
$ws.Range("K33").Value = 1065.2858
$ws.Range("L33").Value = 38
$ws.Range("M33").Value = -836.2858000000001
$ws.Range("N33").Value = -496

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 4273838
$ws.Range("I96").Value = 5128400
$ws.Range("J96").Value = 1029
$ws.Range("K96").Value = 15385200
$ws.Range("L96").Value = 3087
$ws.Range("M96").Value = -15383827
$ws.Range("N96").Value = -5833

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1876
$ws.Range("I100").Value = 1395
$ws.Range("J100").Value = 2196.6667
$ws.Range("K100").Value = 1395
$ws.Range("L100").Value = 2196.6667
$ws.Range("M100").Value = -854
$ws.Range("N100").Value = -3278.6667

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 8578
$ws.Range("I116").Value = 3489.5
$ws.Range("J116").Value = 11485.714
$ws.Range("K116").Value = 3489.5
$ws.Range("L116").Value = 11485.714
$ws.Range("M116").Value = -47.5
$ws.Range("N116").Value = -18369.714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 7126.2324
$ws.Range("I132").Value = 5946.355
$ws.Range("J132").Value = 10174.25
$ws.Range("K132").Value = 17839.065
$ws.Range("L132").Value = 30522.75
$ws.Range("M132").Value = -15309.065
$ws.Range("N132").Value = -35582.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 1040.5946
$ws.Range("I135").Value = 722.7143
$ws.Range("J135").Value = 1234.0869
$ws.Range("K135").Value = 6504.428699999999
$ws.Range("L135").Value = 11106.7821
$ws.Range("M135").Value = -3969.428699999999
$ws.Range("N135").Value = -16176.7821

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 9506.125
$ws.Range("I141").Value = 2914.5386
$ws.Range("K141").Value = 8743.6158
$ws.Range("M141").Value = -3563.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9655.375
$ws.Range("I32").Value = 8667.326999999999
$ws.Range("K32").Value = 8667.326999999999
$ws.Range("M32").Value = -8380.326999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1588
$ws.Range("I102").Value = 1320
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1320
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 302
$ws.Range("N102").Value = -7244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2140.7
$ws.Range("I105").Value = 1602.25
$ws.Range("J105").Value = 2499.6667
$ws.Range("K105").Value = 1602.25
$ws.Range("L105").Value = 2499.6667
$ws.Range("M105").Value = 144.75
$ws.Range("N105").Value = -5993.6667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4977681.5
$ws.Range("I31").Value = 2415.4
$ws.Range("K31").Value = 2415.4
$ws.Range("M31").Value = -2120.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4977681.5
$ws.Range("I34").Value = 2415.4
$ws.Range("K34").Value = 2415.4
$ws.Range("M34").Value = -2213.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5541.1924
$ws.Range("I58").Value = 2854.3845
$ws.Range("J58").Value = 8228
$ws.Range("K58").Value = 2854.3845
$ws.Range("L58").Value = 8228
$ws.Range("M58").Value = -2651.3845
$ws.Range("N58").Value = -8634

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 1272.8572
$ws.Range("I105").Value = 1309.2307
$ws.Range("J105").Value = 800
$ws.Range("K105").Value = 1309.2307
$ws.Range("L105").Value = 800
$ws.Range("M105").Value = 437.7692999999999
$ws.Range("N105").Value = -4294

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1501.6111
$ws.Range("I107").Value = 353.875
$ws.Range("J107").Value = 2419.8
$ws.Range("K107").Value = 353.875
$ws.Range("L107").Value = 2419.8
$ws.Range("M107").Value = 1566.125
$ws.Range("N107").Value = -6259.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2622.6562
$ws.Range("I132").Value = 1831.3125
$ws.Range("J132").Value = 3414
$ws.Range("K132").Value = 5493.9375
$ws.Range("L132").Value = 10242
$ws.Range("M132").Value = -2963.9375
$ws.Range("N132").Value = -15302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 5541.1924
$ws.Range("I136").Value = 2854.3845
$ws.Range("J136").Value = 8228
$ws.Range("K136").Value = 8563.1535
$ws.Range("L136").Value = 24684
$ws.Range("M136").Value = -6013.1535
$ws.Range("N136").Value = -29784

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 58.25
$ws.Range("I14").Value = 58.25
$ws.Range("K14").Value = 174.75
$ws.Range("M14").Value = -1.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 772.5909
$ws.Range("I107").Value = 717.6
$ws.Range("J107").Value = 890.4286
$ws.Range("K107").Value = 717.6
$ws.Range("L107").Value = 890.4286
$ws.Range("M107").Value = 1202.4
$ws.Range("N107").Value = -4730.4286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 10380.083
$ws.Range("J101").Value = 10380.083
$ws.Range("L101").Value = 10380.083
$ws.Range("N101").Value = -16870.083

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H103").Value = 55200.668
$ws.Range("J103").Value = 55200.668
$ws.Range("L103").Value = 55200.668
$ws.Range("N103").Value = -57544.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 14685
$ws.Range("J104").Value = 14685
$ws.Range("L104").Value = 14685
$ws.Range("N104").Value = -21673

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2919.7273
$ws.Range("I136").Value = 1721.4
$ws.Range("J136").Value = 3918.3333
$ws.Range("K136").Value = 5164.200000000001
$ws.Range("L136").Value = 11754.9999
$ws.Range("M136").Value = -2614.200000000001
$ws.Range("N136").Value = -16854.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 672.0625
$ws.Range("I100").Value = 427.27274
$ws.Range("J100").Value = 1210.6
$ws.Range("K100").Value = 854.54548
$ws.Range("L100").Value = 2421.2
$ws.Range("M100").Value = -313.54548
$ws.Range("N100").Value = -3503.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 935
$ws.Range("I113").Value = 220
$ws.Range("J113").Value = 1210
$ws.Range("K113").Value = 660
$ws.Range("L113").Value = 3630
$ws.Range("M113").Value = 1510
$ws.Range("N113").Value = -7970

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 687433.75
$ws.Range("I136").Value = 1278311
$ws.Range("J136").Value = 2016.12
$ws.Range("K136").Value = 3834933
$ws.Range("L136").Value = 6048.36
$ws.Range("M136").Value = -3832383
$ws.Range("N136").Value = -11148.36
